# The workbook gained one new data row: a fresh observation was recorded
# and inserted right where the old row 119 used to be, pushing every
# subsequent row (119-170) down by one (now 120-171). The sheet's used
# range therefore grows from A1:R170 to A1:R171.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 119..170 down to 120..171, leaving row 119 free for the
# newly inserted record.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new observation.
$ws.Range("A119").Value = 9
$ws.Range("B119").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C119").Value = "Metropolitana"
$ws.Range("D119").Value = 44510
$ws.Range("E119").Value = 13
$ws.Range("F119").Value = 100112026
$ws.Range("G119").Value = "Haba"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 52
$ws.Range("K119").Value = 6000
$ws.Range("L119").Value = 7000
$ws.Range("M119").Value = 6500
$ws.Range("N119").Value = "$/saco 25 kilos"
$ws.Range("O119").Value = "Región Metropolitana"
$ws.Range("P119").Value = 260
$ws.Range("Q119").Value = 25
$ws.Range("R119").Value = "Hortaliza"
